# Auto-generated edit script applying the cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to stay a text/string cell (matches original inlineStr type)
    # even when the new value looks like a plain number (e.g. '246.31').
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "35.346.70"
$ws.Range("E2").Value = "  +0.42%  "
Set-TextValue $ws.Range("D3") "1.906.58"
$ws.Range("E4").Value = "  -0.47%  "
Set-TextValue $ws.Range("D5") "246.31"
$ws.Range("E5").Value = "  +3.02%  "
$ws.Range("E6").Value = "  +6.75%  "
$ws.Range("E7").Value = "  -0.47%  "
Set-TextValue $ws.Range("D8") "41.39"
$ws.Range("E8").Value = "  -2.07%  "
Set-TextValue $ws.Range("D9") "0.350"
$ws.Range("E9").Value = "  +6.31%  "
Set-TextValue $ws.Range("D10") "52.73"
$ws.Range("E10").Value = "  +12.38%  "
Set-TextValue $ws.Range("D11") "0.0723"
$ws.Range("E11").Value = "  +4.23%  "
Set-TextValue $ws.Range("D12") "0.0994"
$ws.Range("E12").Value = "  +0.44%  "
Set-TextValue $ws.Range("D13") "2.182.96"
$ws.Range("E13").Value = "  +2.61%  "
Set-TextValue $ws.Range("D14") "12.11"
$ws.Range("E14").Value = "  +5.30%  "
$ws.Range("E15").Value = "  +3.28%  "
Set-TextValue $ws.Range("D16") "1.913.00"
$ws.Range("E16").Value = "  +1.97%  "
$ws.Range("E17").Value = "  +2.92%  "
Set-TextValue $ws.Range("D18") "35.342.53"
$ws.Range("E18").Value = "  +0.55%  "
Set-TextValue $ws.Range("D19") "72.37"
$ws.Range("E19").Value = "  +3.52%  "
Set-TextValue $ws.Range("D20") "0.0₃0823"
$ws.Range("E20").Value = "  +3.45%  "
Set-TextValue $ws.Range("D21") "239.87"
$ws.Range("E21").Value = "  -0.46%  "
Set-TextValue $ws.Range("D22") "12.51"
$ws.Range("E22").Value = "  +2.40%  "
Set-TextValue $ws.Range("D23") "4.84"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("E24").Value = "  -0.48%  "
Set-TextValue $ws.Range("D25") "2.31"
$ws.Range("E25").Value = "  +1.52%  "
Set-TextValue $ws.Range("D26") "2.35"
$ws.Range("E26").Value = "  +23.01%  "
Set-TextValue $ws.Range("D27") "170.57"
$ws.Range("E27").Value = "  +0.69%  "
Set-TextValue $ws.Range("D28") "8.46"
$ws.Range("E28").Value = "  +5.22%  "
Set-TextValue $ws.Range("D29") "18.43"
$ws.Range("E29").Value = "  +4.12%  "
Set-TextValue $ws.Range("D30") "0.128"
$ws.Range("E30").Value = "  +2.92%  "
Set-TextValue $ws.Range("D31") "4.15"
$ws.Range("E31").Value = "  +3.64%  "
Set-TextValue $ws.Range("D32") "0.0566"
$ws.Range("E32").Value = "  +0.97%  "
Set-TextValue $ws.Range("D33") "1.01"
$ws.Range("E33").Value = "  -0.41%  "
Set-TextValue $ws.Range("D34") "0.934"
$ws.Range("E34").Value = "  +14.13%  "
Set-TextValue $ws.Range("D35") "4.12"
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("E36").Value = "  -3.13%  "
Set-TextValue $ws.Range("D37") "2.05"
$ws.Range("E37").Value = "  +0.17%  "
Set-TextValue $ws.Range("D38") "1.34"
$ws.Range("E38").Value = "  +1.94%  "
Set-TextValue $ws.Range("D39") "1.11"
$ws.Range("E39").Value = "  +0.85%  "
$ws.Range("E40").Value = "  +3.30%  "
$ws.Range("E41").Value = "  +8.67%  "
Set-TextValue $ws.Range("D42") "0.0632"
$ws.Range("E42").Value = "  +6.15%  "
Set-TextValue $ws.Range("D43") "90.18"
$ws.Range("E43").Value = "  +0.11%  "
Set-TextValue $ws.Range("D44") "1.342.27"
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E45").Value = "  +2.56%  "
Set-TextValue $ws.Range("D46") "47.59"
$ws.Range("E46").Value = "  +37.44%  "
Set-TextValue $ws.Range("D49") "6.54"
$ws.Range("E49").Value = "  -0.21%  "
Set-TextValue $ws.Range("D50") "2.091.38"
$ws.Range("E50").Value = "  +2.34%  "
Set-TextValue $ws.Range("D51") "0.0707"
$ws.Range("E51").Value = "  +3.37%  "

# Row 47 (previously HuobiToken) now shows the MXToken entry with refreshed figures
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D47") "2.78"
$ws.Range("E47").Value = "  +1.64%  "

# Row 48 (previously MXToken) now shows the HuobiToken entry with refreshed figures
$ws.Range("B48").Value = "HuobiToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D48") "2.40"
$ws.Range("E48").Value = "  -0.96%  "
